$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "66.668.67"
Set-TextValue "E2" "  +2.52%  "
Set-TextValue "D3" "3.204.86"
Set-TextValue "E3" "  +1.46%  "
Set-TextValue "E4" "  +0.08%  "
Set-TextValue "D5" "604.63"
Set-TextValue "E5" "  +4.29%  "
Set-TextValue "D6" "156.91"
Set-TextValue "E6" "  +4.71%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.06%  "
Set-TextValue "D8" "0.558"
Set-TextValue "E8" "  +6.07%  "
Set-TextValue "D9" "3.202.51"
Set-TextValue "E9" "  +1.39%  "
Set-TextValue "E10" "  +1.49%  "
Set-TextValue "D11" "5.85"
Set-TextValue "E11" "  -4.73%  "
Set-TextValue "D12" "0.520"
Set-TextValue "E12" "  +3.86%  "
Set-TextValue "E13" "  +1.97%  "
Set-TextValue "D14" "39.32"
Set-TextValue "E14" "  +5.63%  "
Set-TextValue "D15" "3.734.47"
Set-TextValue "E15" "  +1.60%  "
Set-TextValue "D16" "66.720.43"
Set-TextValue "E16" "  +2.76%  "
Set-TextValue "E17" "  +5.10%  "
Set-TextValue "D18" "3.208.88"
Set-TextValue "E18" "  +1.39%  "
Set-TextValue "D19" "525.98"
Set-TextValue "E19" "  +4.14%  "
Set-TextValue "E20" "  +0.72%  "
Set-TextValue "E21" "  +4.06%  "
Set-TextValue "E22" "  +3.85%  "
Set-TextValue "D23" "8.23"
Set-TextValue "E23" "  +6.34%  "
Set-TextValue "E24" "  -0.89%  "
Set-TextValue "D25" "85.79"
Set-TextValue "E25" "  +1.44%  "
Set-TextValue "E26" "  +0.03%  "
Set-TextValue "D27" "9.30"
Set-TextValue "E27" "  +2.81%  "
Set-TextValue "E28" "  +3.15%  "
Set-TextValue "D29" "2.39"
Set-TextValue "E29" "  +9.14%  "
Set-TextValue "E30" "  +6.66%  "
Set-TextValue "E31" "  +9.11%  "
Set-TextValue "D32" "28.37"
Set-TextValue "E32" "  +2.69%  "
Set-TextValue "D33" "1.24"
Set-TextValue "E33" "  +3.07%  "
Set-TextValue "E34" "  +0.09%  "
Set-TextValue "E35" "  +1.26%  "
Set-TextValue "D36" "513.94"
Set-TextValue "E36" "  +7.25%  "
Set-TextValue "D37" "54.93"
Set-TextValue "E37" "  +0.18%  "
Set-TextValue "E38" "  +1.98%  "
Set-TextValue "D39" "0.0427"
Set-TextValue "E39" "  +2.88%  "
Set-TextValue "E40" "  +8.36%  "
Set-TextValue "D41" "8.92"
Set-TextValue "E41" "  +2.18%  "
Set-TextValue "E42" "  -0.80%  "
Set-TextValue "D43" "0.0₃0683"
Set-TextValue "E43" "  +15.11%  "
Set-TextValue "E44" "  +6.57%  "
Set-TextValue "E45" "  +1.25%  "
Set-TextValue "D46" "2.897.11"
Set-TextValue "E46" "  -3.04%  "
Set-TextValue "D47" "28.68"
Set-TextValue "E47" "  +1.16%  "
Set-TextValue "D48" "2.77"
Set-TextValue "E48" "  +11.57%  "
Set-TextValue "E49" "  +3.66%  "
Set-TextValue "E50" "  -0.04%  "
Set-TextValue "E51" "  +5.04%  "
